# Apply the Team_PER_2012 update:
#  - Team names in column B (rows 2-31) are re-shuffled to a new order
#  - PER values in column C (rows 2-31) are replaced with corrected figures
#
# commit message: "updated code, fixed PER bug and added python ML code"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$teams = @(
    "POR",
    "NJN",
    "CLE",
    "DAL",
    "ATL",
    "OKC",
    "CHA",
    "WAS",
    "MIL",
    "LAC",
    "SAS",
    "DET",
    "ORL",
    "UTA",
    "MEM",
    "HOU",
    "DEN",
    "LAL",
    "GSW",
    "IND",
    "CHI",
    "PHI",
    "BOS",
    "TOR",
    "MIA",
    "SAC",
    "PHO",
    "NOH",
    "NYK",
    "MIN"
)

$perValues = @(
    13.65454545454545,
    11.65,
    11.72857142857143,
    12.6,
    13.40666666666667,
    13.42666666666667,
    10.56,
    10.7125,
    13.21428571428572,
    13.6,
    14.59411764705882,
    12.84285714285714,
    11.22666666666666,
    13.25714285714286,
    11.35333333333333,
    13.5,
    15.35,
    12.20833333333333,
    12.94,
    15.05714285714286,
    15.36428571428572,
    11.89333333333333,
    12.23333333333333,
    12.27333333333333,
    12.36428571428572,
    13.60769230769231,
    14.20769230769231,
    12.63684210526316,
    12.78235294117647,
    13.20666666666667
)

for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $teams[$i]
    $ws.Cells.Item($row, 3).Value = $perValues[$i]
}
